$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# Original: A1=NGUOI MUON, B1=MA SACH, C1=TEN SACH, D1=TAC GIA,
#           E1=NGAY MUON, F1=NGAY HET HAN, G1=TINH TRANG, H1=TIEN PHAT
# New:      A1=CMND, B1=MA SACH, C1=NGAY MUON, D1="NGAY HET HAN ",
#           E1=TINH TRANG, F1 blank
$ws.Range("A1").Value = "CMND"
$ws.Range("C1").Value = "NGAY MUON"
$ws.Range("D1").Value = "NGAY HET HAN "
$ws.Range("E1").Value = "TINH TRANG"
$ws.Range("G1").ClearContents()
$ws.Range("H1").ClearContents()
$ws.Range("F1").NumberFormat = "dd/mm/yyyy"
$ws.Range("F1").Value = ""

# ---- Data rows (2-5) ----
$ws.Range("A2").Value = 216920917
$ws.Range("B2").Value = "KHTN002"
$ws.Range("C2").Value = 43466
$ws.Range("D2").Formula = "=C2+7"
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = 287918095
$ws.Range("B3").Value = "KHTN003"
$ws.Range("C3").Value = 43467
$ws.Range("D3").Formula = "=C3+7"
$ws.Range("E3").Value = 0

$ws.Range("A4").Value = 165024578
$ws.Range("B4").Value = "KHTN008"
$ws.Range("C4").Value = 43468
$ws.Range("D4").Formula = "=C4+7"
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = 275771181
$ws.Range("B5").Value = "KHTN004"
$ws.Range("C5").Value = 43469
$ws.Range("D5").Formula = "=C5+7"
$ws.Range("E5").Value = 1

# ---- Blank placeholder rows (6-20) pre-formatted for new entries ----
for ($r = 6; $r -le 20; $r++) {
    $ws.Range("C$r").NumberFormat = "dd/mm/yyyy"
    $ws.Range("D$r").NumberFormat = "dd/mm/yyyy"
    $ws.Range("E$r").NumberFormat = "0"
}

# ---- Number formats for the real data rows ----
$ws.Range("C2:D5").NumberFormat = "dd/mm/yyyy"
$ws.Range("E2:E5").NumberFormat = "0"

# ---- Column widths (auto-fit to content, like the real edit) ----
$ws.Columns("C:D").AutoFit()

$ws.Range("F5").Select()
